$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, pushing existing rows 14-89 down to 15-90.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new record's data.
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Vega Monumental Concepción"
$ws.Range("C14").Value = "Bíobío"
$ws.Range("D14").Value = (Get-Date -Year 2022 -Month 1 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100108
$ws.Range("H14").Value = "Tropicales y subtropicales"
$ws.Range("I14").Value = 100108002
$ws.Range("J14").Value = "Mango"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 6500
$ws.Range("O14").Value = 7000
$ws.Range("P14").Value = 6750
$ws.Range("Q14").Value = '$/bandeja 4 kilos'
$ws.Range("R14").Value = "Ecuador"
$ws.Range("S14").Value = 1688
$ws.Range("T14").Value = 4
